# Add a new entry (row 15) to the "Journal" work-log sheet and
# update the active cell selection, as captured in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# New row of data: 2024-07-10, 13:00 -> 20:00, category "Analyse et état de l'art"
# (D15 already holds a shared formula C15-B15, so it recalculates automatically)
$ws.Range("A15").Value = 45483
$ws.Range("B15").Value = 0.54166666666666663
$ws.Range("C15").Value = 0.83333333333333337
$ws.Range("E15").Value = "Analyse et état de l'art"

# Move the active cell / selection to H9, matching the saved view state
$ws.Range("H9").Select()
